# Updated cryptos list values (Price column D and Volume(1h) column E)
# applied via headless Excel COM-interop.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Price" values (column D) for the affected rows. These are stored as
# plain text in the workbook (e.g. "26.147.76"), so we force the cell to a
# text number-format before assigning the value and then restore the cell's
# original style, to avoid Excel auto-converting the text into a number.
$priceUpdates = @{
    2 = "26.147.76"
    3 = "1.668.10"
    6 = "0.5216"
    8 = "0.2617"
    9 = "0.06335"
    10 = "21.15"
    11 = "0.07533"
    12 = "1.674.78"
    13 = "4.431"
    14 = "0.5501"
    15 = "66.42"
    16 = "0.000007961"
    17 = "26.156.94"
    19 = "4.721"
    20 = "186.46"
    21 = "10.28"
    22 = "6.191"
    23 = "1.004"
    24 = "149.40"
    25 = "0.1245"
    26 = "7.490"
    28 = "0.06375"
    29 = "1.349"
    30 = "1.274"
    31 = "3.494"
    32 = "3.413"
    33 = "1.639"
    35 = "2.410"
    36 = "0.6019"
    38 = "1.110.55"
    40 = "0.01615"
    41 = "0.8689"
    43 = "99.99"
    44 = "1.819.52"
    45 = "0.00000000107"
    46 = "55.31"
    47 = "1.000"
    48 = "8.028"
    49 = "0.05228"
    50 = "0.4247"
    51 = "5.920"
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Range("D$row")
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
    $cell.Style = $origStyle
}

# New "Volume(1h)" values (column E) for the affected rows. These already
# contain non-numeric characters (spaces and a percent sign) so Excel keeps
# them as plain text without any extra handling.
$volumeUpdates = @{
    2 = "  -0.54%  "
    3 = "  -1.22%  "
    5 = "  -3.66%  "
    6 = "  -2.18%  "
    7 = "  -0.39%  "
    8 = "  -3.64%  "
    9 = "  -1.25%  "
    10 = "  -2.61%  "
    11 = "  -2.09%  "
    12 = "  -1.07%  "
    13 = "  -2.23%  "
    14 = "  -5.07%  "
    15 = "  -0.88%  "
    16 = "  -4.98%  "
    17 = "  -0.70%  "
    18 = "  -0.44%  "
    19 = "  -3.69%  "
    20 = "  -3.68%  "
    21 = "  -5.36%  "
    22 = "  -1.26%  "
    23 = "  -0.40%  "
    24 = "  +0.09%  "
    25 = "  -2.52%  "
    26 = "  -4.60%  "
    27 = "  -0.01%  "
    28 = "  +3.86%  "
    29 = "  -1.96%  "
    30 = "  -3.92%  "
    31 = "  -3.05%  "
    32 = "  -4.67%  "
    33 = "  -3.00%  "
    34 = "  -2.76%  "
    35 = "  -0.72%  "
    36 = "  -2.60%  "
    37 = "  -0.41%  "
    38 = "  +0.02%  "
    39 = "  -2.11%  "
    40 = "  -1.49%  "
    41 = "  -3.19%  "
    42 = "  -0.81%  "
    44 = "  -1.12%  "
    45 = "  -3.04%  "
    46 = "  -4.21%  "
    47 = "  -0.92%  "
    48 = "  -1.17%  "
    50 = "  -1.08%  "
    51 = "  -2.49%  "
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Range("E$row").Value = $volumeUpdates[$row]
}

Write-Host "Updated $($priceUpdates.Count) price cells and $($volumeUpdates.Count) volume cells"
